$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up column L width to match the other data columns ---
$ws.Columns.Item(12).ColumnWidth = 14.29

# --- Base formatting for the whole new column (font, size, centered alignment) ---
# matches the 宋体 / size 12 / center-center formatting used by columns A:K
$fullCol = $ws.Range("L1:L115")
$fullCol.Font.Name = "宋体"
$fullCol.Font.Size = 12
$fullCol.HorizontalAlignment = -4108
$fullCol.VerticalAlignment = -4108

# --- Row 1: new date header "2025/12/01" (kept as text, like the other date headers) ---
$hdr1 = $ws.Range("L1")
$hdr1.NumberFormat = "@"
$hdr1.Value = "2025/12/01"
$hdr1.NumberFormat = "General"

# --- Row 2: index name "上证", bold like the rest of row 2 ---
$hdr2 = $ws.Range("L2")
$hdr2.NumberFormat = "@"
$hdr2.Value = "上证"
$hdr2.NumberFormat = "General"
$hdr2.Font.Bold = $true

# --- Numeric data rows: two-decimal values, same "0.00" number format as column K ---
$numRows = @{
    3 = 61.11;
    4 = 3904.9;
    6 = 47.8;
    7 = 5487.42;
    9 = 52.75;
    10 = 4560.67;
    12 = 55.87;
    13 = 7086.13;
    15 = 26.18;
    16 = 2668.8;
    18 = 96.8;
    19 = 6849.09;
    21 = 70.09999999999999;
    22 = 86034.03999999999;
    24 = 83.29000000000001;
    25 = 23836.79;
    27 = 70.06999999999999;
    28 = 49371.8;
    30 = 47.96;
    31 = 5579.25;
    33 = 9.81;
    34 = 31672.68;
    36 = 27.05;
    37 = 3253.68;
    39 = 41.87;
    40 = 3080.16;
    42 = 15.21;
    43 = 6938.26;
    45 = 27.4;
    46 = 8493.969999999999;
    48 = 11.99;
    49 = 12878.04;
    51 = 24.66;
    52 = 12411.4;
    54 = 21.94;
    55 = 9682.030000000001;
    57 = 26.83;
    58 = 16031.85;
    60 = 31.41;
    61 = 17526.85;
    63 = 21.04;
    64 = 9801.49;
    66 = 14.17;
    67 = 9889.290000000001;
    69 = 23.72;
    70 = 3020.73;
    72 = 41.98;
    73 = 5654.62;
    75 = 25.24;
    76 = 9059.16;
    78 = 13.37;
    79 = 1314.12;
    81 = 54.83;
    82 = 2790.31;
    84 = 58.79;
    85 = 2840.93;
    87 = 50.6;
    88 = 2928.92;
    90 = 43.61;
    91 = 1991.37;
    93 = 27.39;
    94 = 13621.38;
    96 = 84.64;
    97 = 8856.049999999999;
    99 = 56.22;
    100 = 12114.87;
    102 = 5.82;
    103 = 2228.44;
    105 = 25.23;
    106 = 831.99;
    108 = 29.4;
    109 = 2817.43;
    111 = 20.15;
    112 = 3859.53;
    114 = 29.02;
    115 = 3227.35
}
foreach ($r in $numRows.Keys) {
    $cell = $ws.Cells.Item($r, 12)
    $cell.NumberFormat = "0.00"
    $cell.Value = $numRows[$r]
}

# --- Blank separator rows in column L keep the base General format with no value ---
# (already created by the $fullCol formatting pass above; nothing further required)

Write-Host "L column populated"